$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.31%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-5.51%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.833"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-10.03%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05969"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.45%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.689"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.49%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8755"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.05%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'5.22%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1416"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.14%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.03595"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.94%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07203"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.26%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-1.78%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09240"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.06%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001547"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.10%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005959"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.22%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.486"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.28%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.225"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.00%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.218"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-2.36%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "'0.01065"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1,660.13%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3137"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.01%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-0.31%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.523"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.17%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04235"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.56%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1379"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.02%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001221"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.54%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004519"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.10%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.03%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'-22.98%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03850"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.18%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006148"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'7.12%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1105"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.58%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.03%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'6.43%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005489"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.80%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.1091"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'8.95%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002126"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-3.92%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E50").Style = "Normal"
